# Apply the "mansoni_coverage_scenario_2" coverage-update edit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Platform Coverage")
$ws2 = $wb.Worksheets.Item("MarketShare")

# ---------------------------------------------------------------------------
# Sheet 1 "Platform Coverage"
# ---------------------------------------------------------------------------

# Row 2: drop the stray H2 value, bump the rest of the series from 0.6 -> 0.736
$ws1.Range("H2").ClearContents()
foreach ($col in @("J", "L", "N", "P", "R", "T", "V")) {
    $ws1.Range($col + "2").Value = 0.736
}

# Row 3: X3:AZ3 bumped from 0.75 -> 0.92
for ($c = 24; $c -le 52; $c++) {
    $ws1.Cells.Item(3, $c).Value = 0.92
}

# Row 4: X4:AZ4 bumped from 0.5 -> 0.613, new explicit-black-font style applied
for ($c = 24; $c -le 52; $c++) {
    $cell = $ws1.Cells.Item(4, $c)
    $cell.Value = 0.613
    $cell.Font.Color = 0
}

# Row 5: X5:AZ5 bumped from 0.5 -> 0.613, same new style applied
for ($c = 24; $c -le 52; $c++) {
    $cell = $ws1.Cells.Item(5, $c)
    $cell.Value = 0.613
    $cell.Font.Color = 0
}

# Row 9 (new): Vector Control parameters
$ws1.Range("B9").Value = "Vector Control"
$ws1.Range("K9").Value = 0.000000001
$ws1.Range("L9").Value = 0.00000001

# ---------------------------------------------------------------------------
# Sheet 2 "MarketShare"
# ---------------------------------------------------------------------------

# Row 1: switch the yearly header (D:Z, 2018-2040) for a half-yearly one (D:AV, 2018-2040 step 0.5)
$year = 2018
$half = 0
for ($c = 4; $c -le 48; $c++) {
    if ($half -eq 0) {
        $ws2.Cells.Item(1, $c).Value = $year
        $half = 1
    } else {
        $ws2.Cells.Item(1, $c).Value = $year + 0.5
        $year = $year + 1
        $half = 0
    }
}

# Row 3: extend the all-ones series out to the new half-yearly columns (AA3:AV3)
for ($c = 27; $c -le 48; $c++) {
    $ws2.Cells.Item(3, $c).Value = 1
}

# ---------------------------------------------------------------------------
# Window / selection state
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("Y3:AZ3").Select()
$excel.ActiveWindow.Zoom = 120

$ws2.Activate()
$ws2.Range("T3:AV3").Select()
$excel.ActiveWindow.Zoom = 181
